$d = $word.ActiveDocument

# Helper approach: to get Word to emit the inserted/changed text as its own
# run (as real Word does when a user types in the middle of existing text),
# we briefly toggle a character formatting property (Bold on/off) over a
# range before editing its text. Toggling the same value back off leaves the
# run's rPr empty again, but it forces the engine to keep that span as a
# separate run instead of re-coalescing it into its neighbours.

# ----- Edit 1: "It is more than that in it." -> "It is not more than that in it." -----
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("It is more than that in it", $false, $true)
if (-not $found1) {
    throw "Could not find target phrase for edit 1"
}
$phrase1Start = $rng1.Start

# Isolate the whole phrase as its own run first (this also keeps explicit,
# empty rPr elements on the surrounding runs).
$rng1.Bold = 1
$rng1.Bold = 0

# Replace the phrase's text in place with the new wording.
$rng1.Text = "It is not more than that in it"

# Now isolate just the inserted word "not" into its own run.
$not1Start = $phrase1Start + 6
$not1End = $not1Start + 3
$not1Rng = $d.Range($not1Start, $not1End)
$not1Rng.Bold = 1
$not1Rng.Bold = 0

# ----- Edit 2: "we concentrate in pairs" -> "we concentrate on pairs" -----
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("we concentrate in pairs", $false, $true)
if (-not $found2) {
    throw "Could not find target phrase for edit 2"
}
$phrase2Start = $rng2.Start

# Isolate the whole phrase as its own run first.
$rng2.Bold = 1
$rng2.Bold = 0

# Replace the phrase's text in place, swapping "in" for "on".
$rng2.Text = "we concentrate on pairs"

# Now isolate just the new word "on" into its own run.
$on2Start = $phrase2Start + 15
$on2End = $on2Start + 2
$on2Rng = $d.Range($on2Start, $on2End)
$on2Rng.Bold = 1
$on2Rng.Bold = 0

Write-Host "Edits applied"
